# Generate Report for Handoff
# b.md moves from "Handed back: in sync with en-US" to "Ready for handoff"
# after a new (newer) handoff xliff (b.63290e5768f688058c7b37413b0a5c26c308f864.*)
# is generated for it; the previous handback is now stale, so an Error Detail
# note about the stale handback is recorded on the per-locale sheets too.

$wb = $excel.ActiveWorkbook

$statusReadyForHandoff = "Ready for handoff"
$handoffDatetime       = "2016-09-06 10:49:51"

$staleHandbackNote = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/36d1924d66f573b1c2e1316bacdcd4cf13608606/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/46b60f67a7dd6cea1329b24b789d5be55b73875c/e2e/b.md."

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is b.md -- update its zh-cn / de-de status columns and
# the "Latest HO Xliff Generate Date" column.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusReadyForHandoff
$overview.Range("F3").Value = $statusReadyForHandoff
$overview.Range("G3").Value = $handoffDatetime

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 is b.md.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusReadyForHandoff
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-06 10:49:47"
$zhcn.Range("P3").Value = $staleHandbackNote
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------------
# de-de sheet: row 3 is b.md.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusReadyForHandoff
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = $handoffDatetime
$dede.Range("P3").Value = $staleHandbackNote
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
